$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = 'done'
$ws.Range("D29").Value = 'done'
$ws.Range("A31:C31").Clear()
$ws.Range("D30").Value = 'done'
$ws.Range("B33:C33").Clear()
$ws.Range("B34:C34").ClearFormats()
$ws.Range("A34").Value = 'Improve dot fade experience'
$ws.Range("B34").Value = 'body'
$ws.Range("C34").Value = 'high'
$ws.Range("A35").Value = 'dot fade on mobile etc'
$ws.Range("B35").Value = 'body'
$ws.Range("C35").Value = 'high'
$ws.Range("A36").Value = 'full responsive on load'
$ws.Range("B36").Value = 'all'
$ws.Range("C36").Value = 'high'
$ws.Range("B37").Value = 'all'
$ws.Range("C37").Value = 'high'
$ws.Range("B38").Value = 'all'
$ws.Range("C38").Value = 'high'
$ws.Range("A39").Value = 'final checklist'
$ws.Range("B39").Value = 'all'
$ws.Range("C39").Value = 'high'
$ws.Range("E39").Value = 'https://urbanorg.app.box.com/notes/62041029749'
$ws.Range("A30").Value = 'add in the date'
$ws.Range("A40").Value = 'make sure date is correct'
$ws.Range("B40").Value = 'all'
$ws.Range("C40").Value = 'high'
$ws.Range("A32").Value = 'add note to neighborhood chart'
$ws.Range("B32").Value = 'body'
$ws.Range("E32").Value = 'Notes: Urban Institute analysis of student-level data. A student’s neighborhood is defined as the student’s home census tract.'
$ws.Range("A46").Value = 'improve ward map overflow tooltip issue'
$ws.Range("B46").Value = 'maps'
$ws.Range("C46").Value = 'medium'
$ws.Range("A48").Value = 'favicon'
$ws.Range("B48").Value = 'body'
$ws.Range("C48").Value = 'high'
$ws.Range("A57").Value = 'should text under maps say "Note"?'
$ws.Range("B57").Value = 'body'
$ws.Range("C57").Value = 'medium'
$ws.Range("E57").Value = 'Elizabeth'
$ws.Range("A84").Value = 'change distro chart labels to blue'
$ws.Range("A85").Value = 'change distro chart labels to straight and above the thing'
$ws.Range("A86").Value = 'make the axis less crowded on mobile'
$ws.Range("B85").Value = 'mobile'
$ws.Range("A87").Value = 'reduce height on mobile'

$ws.Range("B1").Select()
